$wb = $excel.ActiveWorkbook

# "Overview" sheet: update the Latest HO Xliff Generate Date for row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-30 23:10:39"

# "zh-cn" sheet: update Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2) for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-30 23:10:34"
$wsZhCn.Range("K2").Value = "2016-08-30 23:10:53"

# "de-de" sheet: update Correspond Handoff Datetime (H2, shares the same value as
# Overview!G2) and Correspond Handback DateTime (K2) for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-30 23:10:39"
$wsDeDe.Range("K2").Value = "2016-08-30 23:11:02"
